# H1DR1x-T BOM fixes
# - D5 (Bourns TVS diode) row: fix garbled Description text, add hyperlink on the
#   Octopart URL cell, and grow the row to fit the wrapped description.
# - U2 (MAX14840) row: quantity was blank -> set to 1.
# - R3, R6 row: quantity corrected from 3 to 2.
# - Refresh the view's active cell/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 (Designator D5 / Bourns CDSOD323-T03SC TVS diode) ---------------
# Replace the garbled/duplicated description with the correct one.
$ws.Range("B22").Value = "TVS DIODE 3,3V 10,9V SOD323"

# Row grows to a 2-line wrapped height now that the description is populated.
$ws.Rows.Item(22).RowHeight = 29

# The Octopart URL in E22 becomes a real hyperlink (it was already displaying
# the URL as text but wasn't a clickable hyperlink before).
$ws.Hyperlinks.Add($ws.Range("E22"), "https://octopart.com/cdsod323-t03sc-bourns-10487153?r=sp")

# --- Row 23 (Designator U2 / MAX14840 RS485 transceiver) --------------------
# Quantity was left blank; fill it in with 1.
$ws.Range("F23").Value = 1

# --- Row 14 (Designator R3 , R6) ---------------------------------------------
# Quantity correction: 3 -> 2.
$ws.Range("F14").Value = 2

# --- Refresh selection / scrolled view ---------------------------------------
$ws.Range("D14").Select()
